$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 77

# Columns A and D contain values that look like a date / number
# ("2024-01-19" and "02"), so a leading apostrophe is used to force
# Excel to keep them as literal text (matching the existing text data
# in this column), exactly as a user typing them in would.
$ws.Cells.Item($row, 1).Value = "'2024-01-19"
$ws.Cells.Item($row, 2).Value = "21:49:01"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "'02"
$ws.Cells.Item($row, 5).Value = 138255
$ws.Cells.Item($row, 6).Value = 140566
$ws.Cells.Item($row, 7).Value = 171568
$ws.Cells.Item($row, 8).Value = 148813
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 122570
$ws.Cells.Item($row, 11).Value = 223601
$ws.Cells.Item($row, 12).Value = 255289
$ws.Cells.Item($row, 13).Value = 185371
$ws.Cells.Item($row, 14).Value = 110318
$ws.Cells.Item($row, 15).Value = 41406
$ws.Cells.Item($row, 16).Value = 30924
$ws.Cells.Item($row, 17).Value = 73595
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42887
$ws.Cells.Item($row, 20).Value = -1
